$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlinks that were attached to the old A2/A3 e-mail cells —
# the refreshed data set is plain text, no mailto: links.
$ws.Hyperlinks.Delete()

# A2 used to carry the hyperlink style (blue/underline). Pull a plain
# (unstyled) format from B1 onto A2 before writing the new value, so the
# cell goes back to the regular "Normal" look instead of staying blue.
$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)

# New login rows: column A = e-mail, column B = "password" (shared with B2).
$emails = @(
  "Auttesttt_10@mailinator.com",
  "Auttesttt_11@mailinator.com",
  "Auttesttt_12@mailinator.com",
  "Auttesttt_13@mailinator.com",
  "Auttesttt_14@mailinator.com",
  "Auttesttt_15@mailinator.com",
  "Auttesttt_16@mailinator.com"
)

$row = 2
foreach ($email in $emails) {
  $ws.Cells.Item($row, 1).Value = $email
  $ws.Cells.Item($row, 2).Value = "password"
  $row = $row + 1
}

# Re-select the newly populated block, mirroring the author's selection.
$ws.Range("A3:B14").Select()
